$d = $word.ActiveDocument

$replacements = @(
    @("2025-12-25 Thursday", "2025-12-26 Friday"),
    @("731÷5=146, 1", "909÷8=113, 5"),
    @("301÷7=43, 0", "389÷5=77, 4"),
    @("487÷2=243, 1", "337÷3=112, 1"),
    @("778÷7=111, 1", "348÷8=43, 4"),
    @("948÷6=158, 0", "825÷4=206, 1"),
    @("289÷4=72, 1", "476÷3=158, 2"),
    @("490÷8=61, 2", "108÷8=13, 4"),
    @("564÷4=141, 0", "119÷9=13, 2"),
    @("550÷7=78, 4", "487÷5=97, 2"),
    @("164÷9=18, 2", "175÷8=21, 7"),
    @("872÷2=436, 0", "483÷9=53, 6"),
    @("151÷9=16, 7", "745÷4=186, 1"),
    @("165÷6=27, 3", "885÷4=221, 1"),
    @("225÷3=75, 0", "278÷3=92, 2"),
    @("508÷5=101, 3", "134÷2=67, 0"),
    @("538÷3=179, 1", "137÷5=27, 2"),
    @("668÷4=167, 0", "316÷2=158, 0"),
    @("507÷6=84, 3", "192÷5=38, 2"),
    @("674÷5=134, 4", "682÷3=227, 1"),
    @("222÷7=31, 5", "344÷3=114, 2"),
    @("491÷4=122, 3", "737÷4=184, 1"),
    @("330÷4=82, 2", "186÷3=62, 0"),
    @("722÷9=80, 2", "215÷8=26, 7"),
    @("157÷9=17, 4", "675÷7=96, 3"),
    @("992÷4=248, 0", "820÷3=273, 1")
)

foreach ($pair in $replacements) {
    $old = $pair[0]
    $new = $pair[1]
    $d.Content.Find.Execute($old, $true, $false, $false, $false, $false, $true, 1, $false, $new, 2)
}
